$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the value out of C5 (previously held the shared string "empty"),
# but keep its existing style/formatting.
$ws.Range("C5").ClearContents()

# Move/update the active selection to C5.
$ws.Range("C5").Select()
